# 7.8 History Card & Advanced Story
# Rewrites several dialogue lines on Sheet1 and clears a stray
# "disappear" action flag that had been left in J10.
#
# NOTE: edits are issued in the same order the strings disappear from /
# reappear in the underlying shared-strings table (row order, except
# B15 is applied before B10) so the regenerated sharedStrings.xml keeps
# the exact ordering produced by the original authoring tool.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value  = "Sir, this young lady’s skills seem quite impressive"
$ws.Range("B6").Value  = "Yao, you’re being too hasty."
$ws.Range("B15").Value = "I’m Chen, the best martial artist in this entire manor."
$ws.Range("B10").Value = "I see. You’re amazing——you can spot these details just from a few moves."
$ws.Range("B13").Value = "May I ask your name please?"
$ws.Range("B16").Value = "My apologies, you must be the top disciple of the manor."
$ws.Range("B17").Value = "Top disciple? That’s only because there are just two guards left in Qingliu Manor now."
$ws.Range("B21").Value = "Enough chatting——who are you anyway?"
$ws.Range("B22").Value = "I’m Judge Dee, and this is my student, Yao."
$ws.Range("B32").Value = "Hello!"
$ws.Range("B33").Value = "Hello——may I ask if you’re a physician?"
$ws.Range("B35").Value = "I practice medicine in JiuJiang county at the foot of the mountain. I came up a few days ago for a consultation."
$ws.Range("B36").Value = "Are you familiar with the Lord?"

# Clear the leftover "disappear" marker in J10 (cell keeps its fill style).
$ws.Range("J10").Value = ""

# Match the author's final active selection.
$null = $ws.Range("J10").Select()
